$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for columns that may contain numeric-looking values (B, C, D)
# so Excel does not auto-convert them to numbers, matching the original inlineStr text cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.267.74'

$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.274.83'

$ws.Range("E3").Value = '  +1.43%  '

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '499.75'

$ws.Range("E5").Value = '  +1.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.74'

$ws.Range("E6").Value = '  +1.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'

$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'

$ws.Range("E8").Value = '  -0.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0952'

$ws.Range("E9").Value = '  +0.22%  '

$ws.Range("E10").Value = '  +0.73%  '

$ws.Range("E11").Value = '  +3.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.71'

$ws.Range("E12").Value = '  +1.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.672.57'

$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.68'

$ws.Range("E14").Value = '  +4.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '54.160.45'

$ws.Range("E15").Value = '  +0.17%  '

$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.269.65'

$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.24'

$ws.Range("E18").Value = '  +2.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.14'

$ws.Range("E19").Value = '  +2.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '303.58'

$ws.Range("E20").Value = '  -0.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.33'

$ws.Range("E21").Value = '  -1.82%  '

$ws.Range("E22").Value = '  +0.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.10'

$ws.Range("E23").Value = '  -2.79%  '

$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("E25").Value = '  -0.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.29'

$ws.Range("E26").Value = '  +2.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.70'

$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'PEPE'

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0685'

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Aptos'

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.92'

$ws.Range("E30").Value = '  +0.83%  '

$ws.Range("E31").Value = '  +1.20%  '

$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("E33").Value = '  +0.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.960'

$ws.Range("E34").Value = '  +10.86%  '

$ws.Range("E35").Value = '  +0.68%  '

$ws.Range("E36").Value = '  -1.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.70'

$ws.Range("E37").Value = '  +1.26%  '

$ws.Range("E38").Value = '  -0.90%  '

$ws.Range("E39").Value = '  +0.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.37'

$ws.Range("E40").Value = '  +0.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.83'

$ws.Range("E41").Value = '  -1.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '124.91'

$ws.Range("E42").Value = '  -3.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0492'

$ws.Range("E43").Value = '  +2.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0893'

$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.546'

$ws.Range("E45").Value = '  -0.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '238.69'

$ws.Range("E46").Value = '  -1.20%  '

$ws.Range("E47").Value = '  -0.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0205'

$ws.Range("E48").Value = '  +0.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.75'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.19'

$ws.Range("E50").Value = '  -0.82%  '

$ws.Range("E51").Value = '  -0.40%  '
